# Update master to output generated at c8c62b6
$d = $word.ActiveDocument

# Header date
$d.Content.Find.Execute("2025-10-09 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-10 Friday", 2)

# Row 1 of the division problems table
$d.Content.Find.Execute("432÷2=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "436÷2=", 2)
$d.Content.Find.Execute("455÷7=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "638÷3=", 2)
$d.Content.Find.Execute("735÷9=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "728÷8=", 2)
$d.Content.Find.Execute("494÷5=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "625÷8=", 2)
$d.Content.Find.Execute("540÷7=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "221÷2=", 2)

# Row 2
$d.Content.Find.Execute("671÷3=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "547÷2=", 2)
$d.Content.Find.Execute("752÷5=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "596÷6=", 2)
$d.Content.Find.Execute("957÷2=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "787÷5=", 2)
$d.Content.Find.Execute("740÷7=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "341÷9=", 2)
$d.Content.Find.Execute("781÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "771÷3=", 2)

# Row 3 -- note: "136÷8=" appears twice, so the two cells are addressed
# directly through the table/cell model rather than Find/Replace.
$d.Content.Find.Execute("795÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "957÷9=", 2)
$d.Content.Find.Execute("587÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "955÷9=", 2)
$t = $d.Tables.Item(1)
$t.Cell(9, 3).Range.Text = "302÷4="
$t.Cell(9, 4).Range.Text = "583÷4="
$d.Content.Find.Execute("104÷5=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "388÷9=", 2)

# Row 4
$d.Content.Find.Execute("517÷9=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "499÷8=", 2)
$d.Content.Find.Execute("268÷8=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "440÷5=", 2)
$d.Content.Find.Execute("530÷9=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "647÷9=", 2)
$d.Content.Find.Execute("918÷6=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "695÷9=", 2)
$d.Content.Find.Execute("240÷2=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "470÷6=", 2)

# Row 5
$d.Content.Find.Execute("150÷3=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "982÷3=", 2)
$d.Content.Find.Execute("327÷5=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "283÷2=", 2)
$d.Content.Find.Execute("242÷7=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "463÷5=", 2)
$d.Content.Find.Execute("581÷6=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "372÷9=", 2)
$d.Content.Find.Execute("160÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "788÷5=", 2)
